# Add a new worksheet "total_staff" after the existing "annual_budget" sheet,
# populate it with header/data rows, and leave it as the active/selected sheet
# (matching the diff against supervision_metrics.xlsx).

$wb = $excel.ActiveWorkbook

# Insert the new sheet at the end of the workbook (after the last existing sheet)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "total_staff"

# Header row
$newSheet.Range("A1").Value = "year"
$newSheet.Range("B1").Value = "system"
$newSheet.Range("C1").Value = "value"
$newSheet.Range("D1").Value = "staff_type"

# Data row 2
$newSheet.Range("A2").Value = 2021
$newSheet.Range("B2").Value = "both"
$newSheet.Range("C2").Value = 100

# Data row 3
$newSheet.Range("A3").Value = 2021
$newSheet.Range("B3").Value = "both"
$newSheet.Range("C3").Value = 50

# Populate the staff_type column last, D3 ("SUPPORT") before D2 ("SUPERVISION"),
# so new shared-string entries register in the same order as the source edit.
$newSheet.Range("D3").Value = "SUPPORT"
$newSheet.Range("D2").Value = "SUPERVISION"

# Leave selection on D2 / total_staff as the active sheet & tab
[void]$newSheet.Range("D2").Select()
